$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column header C1, merging the "pot"/"udbr" split columns back
# into a single "TotGoednabDyr_kt_år" column.
$ws.Range("C1").Value = "TotGoednabDyr_kt_år"

# Remove column D entirely (the "...udbr" data introduced by the
# reverted commit), shifting any cells to the right of it left.
$ws.Columns("D").Delete()
